# 4.2 Hardware Interfaces — remove "Draft 2" block entirely, keeping only
# the reviewed "Draft 1" content (commit: "Deleted draft (2), Reviewed
# chosen draft (1) for (4.2 Hardware Interfaces)").
#
# The document currently contains, in order:
#   Draft 1: heading + "4.2 Hardware Interfaces" heading + 3 body paragraphs
#   (blank paragraph)
#   Draft 2: heading + "4.2 Hardware Interfaces" heading + 3 body paragraphs
#   (two trailing blank paragraphs)
#
# Everything from the blank paragraph that follows Draft 1's last sentence
# ("...or backend processing.") through to the end of the document body is
# the Draft 2 block (including its own trailing blanks) and must go.

$d = $word.ActiveDocument

# Anchor on the tail of Draft 1's last sentence so this does not depend on
# hard-coded paragraph indices.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("or backend processing.", $false, $false, $false, $false, `
                      $false, $true, 1, $false, "", 0) | Out-Null

# Move past the paragraph mark that ends that sentence's paragraph, so the
# deletion starts at the very next paragraph (the blank line preceding the
# "Draft 2" heading) ...
$anchor.Collapse(0)            # wdCollapseEnd
$anchor.MoveEnd(1, 1) | Out-Null   # wdCharacter: swallow the paragraph mark
$deleteStart = $anchor.End

# ... and run through to the end of the document's main story (the final
# sectPr is not part of Content, so this removes every remaining paragraph,
# including the two trailing blanks after Draft 2).
$deleteEnd = $d.Content.End

$d.Range($deleteStart, $deleteEnd).Delete()
